# Weekly price-sheet update: a new (most recent) price record is inserted
# at row 17 (the data is kept sorted, newest first, after the header row),
# pushing the previously-existing rows 17-52 down to rows 18-53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 17, shifting rows 17:52 -> 18:53
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with the new weekly record.
# Mercado ID / Mercado / Región / Codreg / Categoría ID / Categoría / Variedad /
# Calidad / prices / unit / origin / clasificación repeat the same
# "Perú - $/malla 20 kilos - 18000" pattern seen throughout this sheet;
# only the date (D) and the volume (J) are new for this entry.
$ws.Cells.Item(17, 1).Value = 10
$ws.Cells.Item(17, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(17, 3).Value = "La Araucanía"
$ws.Cells.Item(17, 4).Value = 44620
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = 100114002
$ws.Cells.Item(17, 7).Value = "Camote"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 50
$ws.Cells.Item(17, 11).Value = 18000
$ws.Cells.Item(17, 12).Value = 18000
$ws.Cells.Item(17, 13).Value = 18000
$ws.Cells.Item(17, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(17, 15).Value = "Perú"
$ws.Cells.Item(17, 16).Value = 900
$ws.Cells.Item(17, 17).Value = 20
$ws.Cells.Item(17, 18).Value = "Hortaliza"
